$p = $ppt.ActivePresentation

# Locate the title-slide subtitle shape that currently reads "FME 2016 Training"
# and update it for the 2017 training edition, matching the target's two-run
# split ("FME " + "2017") produced when PowerPoint edits existing text in place.
$targetText = "FME 2016 Training"
$found = $false

for ($si = 1; $si -le $p.Slides.Count -and -not $found; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count -and -not $found; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.HasText) {
                if ($shape.TextFrame.TextRange.Text -eq $targetText) {
                    $tr = $shape.TextFrame.TextRange
                    $tr.Text = "FME "
                    [void]$tr.InsertAfter("2017")
                    $found = $true
                }
            }
        }
    }
}
